# Update "Number of Mitophagy Events Outside Dialated Area" (column E) values
# for several rows as part of redoing the analysis script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 7
$ws.Range("E7").Value = 4
$ws.Range("E39").Value = 2
$ws.Range("E75").Value = 7
$ws.Range("E79").Value = 2
$ws.Range("E99").Value = 1
$ws.Range("E105").Value = 6
$ws.Range("E111").Value = 3
$ws.Range("E121").Value = 2
$ws.Range("E123").Value = 6
$ws.Range("E131").Value = 4
